# Auto-generated edit script applying the committed cell-value changes
# to the Gilgamesh_Profits workbook (8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H92").Value = 373.1111
$ws.Range("I92").Value = 388.4375
$ws.Range("K92").Value = 388.4375
$ws.Range("M92").Value = 859.5625
$ws.Range("H131").Value = 771521.75
$ws.Range("I131").Value = 1430139.9
$ws.Range("J131").Value = 3134
$ws.Range("K131").Value = 4290419.699999999
$ws.Range("L131").Value = 9402
$ws.Range("M131").Value = -4285379.699999999
$ws.Range("N131").Value = -19482
$ws.Range("H132").Value = 8348.368
$ws.Range("I132").Value = 8706.611000000001
$ws.Range("K132").Value = 26119.833
$ws.Range("M132").Value = -23589.833
$ws.Range("H135").Value = 1673.4546
$ws.Range("I135").Value = 1521.6666
$ws.Range("J135").Value = 1855.6
$ws.Range("K135").Value = 13694.9994
$ws.Range("L135").Value = 16700.4
$ws.Range("M135").Value = -11159.9994
$ws.Range("N135").Value = -21770.4
$ws.Range("H138").Value = 2946.6606
$ws.Range("I138").Value = 2655.4
$ws.Range("J138").Value = 3108.4722
$ws.Range("K138").Value = 7966.200000000001
$ws.Range("L138").Value = 9325.4166
$ws.Range("M138").Value = -2826.200000000001
$ws.Range("N138").Value = -19605.4166
$ws.Range("H141").Value = 2075.1333
$ws.Range("I141").Value = 1945.5454
$ws.Range("J141").Value = 2431.5
$ws.Range("K141").Value = 5836.6362
$ws.Range("L141").Value = 7294.5
$ws.Range("M141").Value = -656.6361999999999
$ws.Range("N141").Value = -17654.5

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2253.6956
$ws.Range("I32").Value = 2266.2444
$ws.Range("J32").Value = 1689
$ws.Range("K32").Value = 2266.2444
$ws.Range("L32").Value = 1689
$ws.Range("M32").Value = -1979.2444
$ws.Range("N32").Value = -2263
$ws.Range("H45").Value = 22541.6
$ws.Range("I45").Value = 26194.295
$ws.Range("K45").Value = 26194.295
$ws.Range("M45").Value = -25817.295
$ws.Range("H60").Value = 32525.25
$ws.Range("I60").Value = 33367
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 33367
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -32634
$ws.Range("N60").Value = -31466
$ws.Range("H61").Value = 2323.2058
$ws.Range("I61").Value = 1263.4546
$ws.Range("K61").Value = 1263.4546
$ws.Range("M61").Value = -1051.4546
$ws.Range("H74").Value = 180549.73
$ws.Range("I74").Value = 253584.77
$ws.Range("J74").Value = 2019.6666
$ws.Range("K74").Value = 253584.77
$ws.Range("L74").Value = 2019.6666
$ws.Range("M74").Value = -252710.77
$ws.Range("N74").Value = -3767.6666
$ws.Range("H77").Value = 180549.73
$ws.Range("I77").Value = 253584.77
$ws.Range("J77").Value = 2019.6666
$ws.Range("K77").Value = 1267923.85
$ws.Range("L77").Value = 10098.333
$ws.Range("M77").Value = -1263555.85
$ws.Range("N77").Value = -18834.333
$ws.Range("H102").Value = 3517.6924
$ws.Range("J102").Value = 5999
$ws.Range("L102").Value = 5999
$ws.Range("N102").Value = -9243
$ws.Range("H110").Value = 482.5
$ws.Range("I110").Value = 482.5
$ws.Range("K110").Value = 482.5
$ws.Range("M110").Value = 1562.5
$ws.Range("H132").Value = 1876.9166
$ws.Range("I132").Value = 1646.125
$ws.Range("K132").Value = 4938.375
$ws.Range("M132").Value = -2408.375
$ws.Range("H136").Value = 2323.2058
$ws.Range("I136").Value = 1263.4546
$ws.Range("K136").Value = 3790.3638
$ws.Range("M136").Value = -1240.3638

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H64").Value = 851.55554
$ws.Range("I64").Value = 429.75
$ws.Range("K64").Value = 429.75
$ws.Range("M64").Value = -204.75
$ws.Range("H67").Value = 851.55554
$ws.Range("I67").Value = 429.75
$ws.Range("K67").Value = 429.75
$ws.Range("M67").Value = 350.25
$ws.Range("H134").Value = 6690.5
$ws.Range("I134").Value = 2544.2307
$ws.Range("K134").Value = 7632.6921
$ws.Range("M134").Value = -5097.6921

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 4466.1704
$ws.Range("I31").Value = 5377.6
$ws.Range("J31").Value = 4219.838
$ws.Range("K31").Value = 5377.6
$ws.Range("L31").Value = 4219.838
$ws.Range("M31").Value = -5082.6
$ws.Range("N31").Value = -4809.838
$ws.Range("H34").Value = 4466.1704
$ws.Range("I34").Value = 5377.6
$ws.Range("J34").Value = 4219.838
$ws.Range("K34").Value = 5377.6
$ws.Range("L34").Value = 4219.838
$ws.Range("M34").Value = -5175.6
$ws.Range("N34").Value = -4623.838
$ws.Range("H58").Value = 2226.625
$ws.Range("I58").Value = 1728
$ws.Range("K58").Value = 1728
$ws.Range("M58").Value = -1525
$ws.Range("H132").Value = 5654858
$ws.Range("I132").Value = 4661.844
$ws.Range("J132").Value = 23816204
$ws.Range("K132").Value = 13985.532
$ws.Range("L132").Value = 71448612
$ws.Range("M132").Value = -11455.532
$ws.Range("N132").Value = -71453672
$ws.Range("H134").Value = 3142.7441
$ws.Range("I134").Value = 2803.5938
$ws.Range("K134").Value = 8410.7814
$ws.Range("M134").Value = -5875.7814
$ws.Range("H136").Value = 2226.625
$ws.Range("I136").Value = 1728
$ws.Range("K136").Value = 5184
$ws.Range("M136").Value = -2634

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 1665.8334
$ws.Range("I132").Value = 968
$ws.Range("J132").Value = 1898.4445
$ws.Range("K132").Value = 8712
$ws.Range("L132").Value = 17086.0005
$ws.Range("M132").Value = -6182
$ws.Range("N132").Value = -22146.0005

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H43").Value = 6995.3335
$ws.Range("I43").Value = 6995.3335
$ws.Range("K43").Value = 6995.3335
$ws.Range("M43").Value = -6844.3335
$ws.Range("H109").Value = 67499.5
$ws.Range("J109").Value = 67499.5
$ws.Range("L109").Value = 67499.5
$ws.Range("N109").Value = -69579.5
$ws.Range("H113").Value = 5759.25
$ws.Range("I113").Value = 3755.5715
$ws.Range("J113").Value = 19785
$ws.Range("K113").Value = 3755.5715
$ws.Range("L113").Value = 19785
$ws.Range("M113").Value = -1585.5715
$ws.Range("N113").Value = -24125
$ws.Range("H132").Value = 2405.2917
$ws.Range("J132").Value = 2802.8
$ws.Range("L132").Value = 8408.400000000001
$ws.Range("N132").Value = -13468.4

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 86957.25
$ws.Range("I40").Value = 94407.91
$ws.Range("K40").Value = 94407.91
$ws.Range("M40").Value = -94271.91
$ws.Range("H55").Value = 1194.4615
$ws.Range("I55").Value = 1014.3333
$ws.Range("J55").Value = 1348.8572
$ws.Range("K55").Value = 1014.3333
$ws.Range("L55").Value = 1348.8572
$ws.Range("M55").Value = -841.3333
$ws.Range("N55").Value = -1694.8572
$ws.Range("H61").Value = 3321.5293
$ws.Range("I61").Value = 3337.5334
$ws.Range("J61").Value = 3201.5
$ws.Range("K61").Value = 3337.5334
$ws.Range("L61").Value = 3201.5
$ws.Range("M61").Value = -3135.5334
$ws.Range("N61").Value = -3605.5
$ws.Range("H94").Value = 29666.334
$ws.Range("J94").Value = 29666.334
$ws.Range("L94").Value = 29666.334
$ws.Range("N94").Value = -31018.334
$ws.Range("H113").Value = 3321.5293
$ws.Range("I113").Value = 3337.5334
$ws.Range("J113").Value = 3201.5
$ws.Range("K113").Value = 3337.5334
$ws.Range("L113").Value = 3201.5
$ws.Range("M113").Value = -1167.5334
$ws.Range("N113").Value = -7541.5
$ws.Range("H122").Value = 9664.4375
$ws.Range("I122").Value = 5176.7144
$ws.Range("K122").Value = 15530.1432
$ws.Range("M122").Value = -13080.1432
$ws.Range("H132").Value = 6710.231
$ws.Range("I132").Value = 8400.6
$ws.Range("J132").Value = 5653.75
$ws.Range("K132").Value = 25201.8
$ws.Range("L132").Value = 16961.25
$ws.Range("M132").Value = -22671.8
$ws.Range("N132").Value = -22021.25
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5800
$ws.Range("K136").Value = 17400
$ws.Range("M136").Value = -14850

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H30").Value = 30009
$ws.Range("I30").Value = 30009
$ws.Range("K30").Value = 30009
$ws.Range("M30").Value = -29902
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H81").Value = 6974.25
$ws.Range("I81").Value = 6974.25
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 13948.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -12887.5
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 6974.25
$ws.Range("I84").Value = 6974.25
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 69742.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -64438.5
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 1105.5625
$ws.Range("I107").Value = 966.55554
$ws.Range("J107").Value = 1284.2858
$ws.Range("K107").Value = 2899.66662
$ws.Range("L107").Value = 3852.8574
$ws.Range("M107").Value = -979.66662
$ws.Range("N107").Value = -7692.857400000001
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 5447.857
$ws.Range("I132").Value = 6031.4165
$ws.Range("J132").Value = 1946.5
$ws.Range("K132").Value = 18094.2495
$ws.Range("L132").Value = 5839.5
$ws.Range("M132").Value = -15564.2495
$ws.Range("N132").Value = -10899.5
$ws.Range("H136").Value = 254280.78
$ws.Range("J136").Value = 915560.9
$ws.Range("L136").Value = 2746682.7
$ws.Range("N136").Value = -2751782.7

